# Auto-generated edit script: updates cryptos price/volume columns
# per the commit "Updated cryptos list on Sun Feb 19 20:56:07 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.634.48"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "'1.691.80"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'314.04"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "'0.3898"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("D8").Value = "'0.4020"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").Value = "'1.495"
$ws.Range("E9").Value = "  +0.95%  "

$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").Value = "'52.86"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "'0.08718"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "'7.612"
$ws.Range("E13").Value = "  +5.22%  "

$ws.Range("D14").Value = "'24.78"
$ws.Range("E14").Value = "  +6.11%  "

$ws.Range("D15").Value = "'7.955"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "'0.00001337"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").Value = "'1.666.37"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").Value = "'98.07"
$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").Value = "'0.07113"
$ws.Range("E19").Value = "  +1.56%  "

$ws.Range("D20").Value = "'19.68"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "'7.235"
$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'14.16"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("D24").Value = "'24.591.84"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'2.996"
$ws.Range("E25").Value = "  -8.94%  "

$ws.Range("D26").Value = "'2.349"
$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("D27").Value = "'22.57"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "'161.19"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").Value = "'8.481"
$ws.Range("E29").Value = "  +12.48%  "

$ws.Range("D30").Value = "'5.226"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'136.13"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").Value = "'1.858.66"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("D33").Value = "'0.08768"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("D34").Value = "'7.399"
$ws.Range("E34").Value = "  +3.78%  "

$ws.Range("D35").Value = "'1.033"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").Value = "'1.985"
$ws.Range("E36").Value = "  +5.85%  "

$ws.Range("D37").Value = "'0.02888"
$ws.Range("E37").Value = "  +6.84%  "

$ws.Range("D38").Value = "'0.2710"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("D39").Value = "'10.73"
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").Value = "'0.09097"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").Value = "'14.05"
$ws.Range("E41").Value = "  -1.99%  "

$ws.Range("D42").Value = "'0.7750"
$ws.Range("E42").Value = "  +2.34%  "

$ws.Range("D43").Value = "'1.454"
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").Value = "'16.47"
$ws.Range("E44").Value = "  +3.36%  "

$ws.Range("D45").Value = "'0.7139"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("D49").Value = "'1.343"
$ws.Range("E49").Value = "  +2.60%  "

$ws.Range("D50").Value = "'137.74"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").Value = "'90.64"
$ws.Range("E51").Value = "  +1.58%  "

